$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"

$ws.Range("D2").Value = "245.64"

$ws.Range("D4").Value = "5.399"

$ws.Range("D5").Value = "0.05816"

$ws.Range("D6").Value = "3.377"

$ws.Range("D7").Value = "6.335"

$ws.Range("D8").Value = "0.8072"

$ws.Range("D9").Value = "0.9772"

$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "0.1426"
$ws.Range("E10").Value = "9WazirXWRX"

$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "0.07491"
$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"

$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D12").Value = "0.03220"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"

$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "0.03033"
$ws.Range("E13").Value = "12BitrueCoinBTR"

$ws.Range("B14").Value = "MCDex"
$ws.Range("C14").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D14").Value = "4.168"
$ws.Range("E14").Value = "13MCDexMCB"

$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").Value = "0.09398"
$ws.Range("E15").Value = "14BitMartTokenBMX"

$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").Value = "0.001597"
$ws.Range("E16").Value = "15BitForexTokenBF"

$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D17").Value = "0.04806"
$ws.Range("E17").Value = "16CoinExTokenCET"

$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D18").Value = "0.0005890"
$ws.Range("E18").Value = "17OneONEWorstin24h"

$ws.Range("D19").Value = "0.006250"

$ws.Range("D20").Value = "0.004090"
$ws.Range("E20").Value = "19HotbitTokenHTB"

$ws.Range("D21").Value = "0.0009955"

$ws.Range("D27").Value = "0.0003593"
$ws.Range("E27").Value = "26UpBotsUBXTBestin24h"

$ws.Range("D40").Value = "0.03890"

$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").Value = "0.1073"
$ws.Range("E41").Value = "40BKEXTokenBKK"

$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").Value = "0.002640"
$ws.Range("E42").Value = "41CEJICEJI"

$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D43").Value = "0.003039"
$ws.Range("E43").Value = "42KickTokenKICK"

$ws.Range("D44").Value = "0.006391"

$ws.Range("D45").Value = "0.00005592"

$ws.Range("D48").Value = "0.1453"

$ws.Range("D2").Style = "Normal"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D48").Style = "Normal"